$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped from
# 45224 (2023-10-25) to 45233 (2023-11-03) for every data row (rows 2-120).
$oldValue = 45224
$newValue = 45233

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 120 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
